$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new blank row above row 39 (pushes the "Questions" block down by one row)
$ws.Rows("39:39").Insert()

# Row 38 (formerly blank) now holds the new "addAnswerOption" localization entry
$ws.Cells.Item(38, 2).Value = "addAnswerOption"
$ws.Cells.Item(38, 3).Value = "Click to add answer option"

# Update the visible selection to match the saved view state
$ws.Range("D30").Select()
